# export client pertama push
# Rename the single worksheet from "Sheet1" to "Adira Finance", drop the
# stale external workbook link (and the "HARI" defined name that pointed
# into it), and let the sheet-scoped Print_Area / _FilterDatabase names
# pick up the new sheet name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Adira Finance"

# Remove the defined name that referenced the external workbook ("HARI").
$names = $wb.Names
for ($i = $names.Count; $i -ge 1; $i--) {
    $n = $names.Item($i)
    if ($n.Name -eq "HARI") {
        $n.Delete()
    }
}

# Break the external link to the old LEMBUR5.xls workbook entirely, which
# also drops the <externalReferences>/externalLink1.xml parts on save.
$links = $wb.LinkSources(1)
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}

# Re-point the Print_Area name at the renamed sheet (_FilterDatabase already
# followed the rename automatically).
$ws.PageSetup.PrintArea = "`$A`$1:`$AA`$5"
